$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "310.61"
Set-TextValue $ws.Range("E2") "2.15%"
Set-TextValue $ws.Range("D3") "38.69"
Set-TextValue $ws.Range("E3") "8.17%"
Set-TextValue $ws.Range("D4") "5.126"
Set-TextValue $ws.Range("E4") "1.77%"
Set-TextValue $ws.Range("D5") "0.08174"
Set-TextValue $ws.Range("E5") "2.64%"
Set-TextValue $ws.Range("D6") "2.009"
Set-TextValue $ws.Range("E6") "8.26%"
Set-TextValue $ws.Range("D7") "7.924"
Set-TextValue $ws.Range("E7") "2.14%"
Set-TextValue $ws.Range("D8") "0.9321"
Set-TextValue $ws.Range("E8") "1.24%"
Set-TextValue $ws.Range("D9") "0.1404"
Set-TextValue $ws.Range("E9") "4.53%"
Set-TextValue $ws.Range("D10") "0.1953"
Set-TextValue $ws.Range("E10") "3.43%"
Set-TextValue $ws.Range("D11") "0.09113"
Set-TextValue $ws.Range("E11") "0.35%"
Set-TextValue $ws.Range("D12") "0.03459"
Set-TextValue $ws.Range("E12") "0.64%"
Set-TextValue $ws.Range("D13") "0.09857"
Set-TextValue $ws.Range("E13") "0.06%"
Set-TextValue $ws.Range("D14") "0.001417"
Set-TextValue $ws.Range("E14") "0.68%"
Set-TextValue $ws.Range("D15") "0.005996"
Set-TextValue $ws.Range("E15") "-1.10%"
Set-TextValue $ws.Range("D16") "3.570"
Set-TextValue $ws.Range("E16") "-4.59%"
Set-TextValue $ws.Range("D17") "4.197"
Set-TextValue $ws.Range("E17") "2.05%"
Set-TextValue $ws.Range("D18") "3.440"
Set-TextValue $ws.Range("E18") "1.34%"
Set-TextValue $ws.Range("D19") "0.3451"
Set-TextValue $ws.Range("E19") "0.24%"
Set-TextValue $ws.Range("D20") "0.1315"
Set-TextValue $ws.Range("E20") "-1.81%"
Set-TextValue $ws.Range("D21") "4.808"
Set-TextValue $ws.Range("E21") "-6.86%"
Set-TextValue $ws.Range("D22") "0.2457"
Set-TextValue $ws.Range("E22") "4.62%"
Set-TextValue $ws.Range("D23") "0.04466"
Set-TextValue $ws.Range("E23") "1.00%"
Set-TextValue $ws.Range("D24") "0.001239"
Set-TextValue $ws.Range("E24") "1.84%"
Set-TextValue $ws.Range("E25") "-9.81%"
Set-TextValue $ws.Range("D27") "0.0001302"
Set-TextValue $ws.Range("E27") "0.22%"
Set-TextValue $ws.Range("D39") "0.02121"
Set-TextValue $ws.Range("E39") "9.20%"
Set-TextValue $ws.Range("D40") "0.05184"
Set-TextValue $ws.Range("E40") "-3.37%"
Set-TextValue $ws.Range("D41") "0.007475"
Set-TextValue $ws.Range("E41") "-1.78%"
Set-TextValue $ws.Range("E42") "-0.80%"
Set-TextValue $ws.Range("E43") "0.92%"
Set-TextValue $ws.Range("D44") "0.002133"
Set-TextValue $ws.Range("E44") "-0.26%"
Set-TextValue $ws.Range("D45") "0.009766"
Set-TextValue $ws.Range("E45") "-4.05%"
Set-TextValue $ws.Range("D46") "0.00006339"
Set-TextValue $ws.Range("E46") "3.02%"
Set-TextValue $ws.Range("D47") "0.00000000751"
Set-TextValue $ws.Range("E47") "0.23%"
Set-TextValue $ws.Range("E48") "-0.24%"
Set-TextValue $ws.Range("D49") "0.001604"
Set-TextValue $ws.Range("E49") "-3.17%"
Set-TextValue $ws.Range("D50") "0.00002103"
Set-TextValue $ws.Range("E50") "0.23%"
Set-TextValue $ws.Range("D51") "0.0002003"
Set-TextValue $ws.Range("E51") "0.23%"
